$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing rawFastqDir / qualityFastqDir values to the new project paths ---
$ws.Cells.Item(2, 2).Value  = "/Users/ruben/Dropbox/Projects/RNAseqPipeline/Data/"
$ws.Cells.Item(3, 2).Value  = "/Users/ruben/Dropbox/Projects/RNAseqPipeline/Results/"
$ws.Cells.Item(12, 2).Value = "/Users/ruben/Dropbox/Projects/RNAseqPipeline/Data/"

# --- Add the new DESeq2 section (rows 30-33) ---
$ws.Cells.Item(30, 1).Value = "DESeq2"
$ws.Cells.Item(30, 1).Font.Bold = $true

$ws.Cells.Item(31, 1).Value = "countTableOrigin"
$ws.Cells.Item(31, 2).Value = "featureCounts"
$ws.Cells.Item(31, 3).Value = "input count table for DESeq2 : featureCounts or HTSeq "

$ws.Cells.Item(32, 1).Value = "CtrlTag"
$ws.Cells.Item(32, 2).Value = "Ctrl"
$ws.Cells.Item(32, 3).Value = "unique tag in name of control samples"

$ws.Cells.Item(33, 1).Value = "TreatmentTag"
$ws.Cells.Item(33, 2).Value = "KO"
$ws.Cells.Item(33, 3).Value = "unique tag in name of treatment samples"

# --- Update the selection / active cell ---
$null = $ws.Range("B15").Select()
